# Append 3 new daily rows (227-229) to the Cavezzo report sheet, mirroring
# the existing column layout: A=date serial (styled like A226), B/C=counts,
# D=computed rate value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 227; A = 44301; B = 4; C = 9;  D = 128.022759601707 },
    @{ Row = 228; A = 44302; B = 2; C = 10; D = 142.2475106685633 },
    @{ Row = 229; A = 44303; B = 0; C = 10; D = 142.2475106685633 }
)

# Copy the format of the last existing date cell (A226) once; reuse it for
# each new A-cell in the appended rows so the row style (borders, bold,
# centered alignment, yyyy-mm-dd date format) matches the prior rows exactly.
$ws.Range("A226").Copy()

foreach ($r in $newRows) {
    $aCell = $ws.Cells.Item($r.Row, 1)
    $aCell.Value = $r.A
    $aCell.PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
